# This script converts the merge-field placeholders that were previously
# split across "${" / "<name>" / "}" runs (with spell-check <w:proofErr/>
# wrappers around the field name) into their tidy final form.
#
#   ${departmentFull}  ->  ${support}     (field renamed, proofErr removed)
#   ${ilos}            ->  ${ilos}        (runs merged, proofErr removed)
#   ${budgetSource}    ->  ${budgetSource}(runs merged, proofErr removed)
#   ${sig_cscp}         -> ${sig_cscp}     (runs merged, proofErr removed)
#   ${sig_csca}         -> ${sig_csca}     (runs merged, proofErr removed)
#   ${sig_sscp}         -> ${sig_sscp}     (runs merged, proofErr removed)
#   ${sig_dean}         -> ${sig_dean}     (runs merged, proofErr removed)
#
# NOTE: double-quoted PowerShell strings interpolate "${name}" as a
# variable reference, so every literal "${...}" below uses single quotes.

$d = $word.ActiveDocument

# --- 1. "departmentFull" field is renamed to "support" -----------------
# Collapse "${departmentFull}" down to "${support}" first (this merges the
# three runs into one clean run and drops the spell-check proofErr tags
# around the old field name), then re-split "${" / "support" / "}" back
# into three separate runs -- matching the original run layout -- by
# nudging (and restoring) a character formatting property on just the
# "support" word so Word has to break it into its own run again.
$rFull = $d.Content
$rFull.Find.Execute('${departmentFull}', $true, $false, $false, $false, $false, $true, 1, $false, '${support}', 2) | Out-Null

$rWord = $d.Content
$rWord.Find.Execute('support', $true, $false, $false, $false, $false, $true, 1, $false, '', 0) | Out-Null
$rWord.Font.Bold = 1
$rWord.Font.Bold = 0

# --- 2. Simple "${<name>}" placeholders: merge the split runs ----------
# Each of these already reads correctly as text; replacing the full
# "${name}" span with itself merges the "${" / name / "}" runs that were
# wrapped in <w:proofErr/> spell-check markers into a single clean run.
$targets = '${ilos}', '${budgetSource}', '${sig_cscp}', '${sig_csca}', '${sig_sscp}', '${sig_dean}'
foreach ($t in $targets) {
    $r = $d.Content
    $r.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2) | Out-Null
}
